$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as text values,
# preserving the original text-based (non-numeric) cell representation.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '63.892.21'
Set-TextValue 'E2' '  +0.19%  '
Set-TextValue 'D3' '3.142.11'
Set-TextValue 'E3' '  +0.77%  '
Set-TextValue 'E4' '  +0.22%  '
Set-TextValue 'D5' '590.87'
Set-TextValue 'E5' '  +0.67%  '
Set-TextValue 'D6' '145.30'
Set-TextValue 'E6' '  -0.81%  '
Set-TextValue 'E7' '  +0.06%  '
Set-TextValue 'D8' '3.133.23'
Set-TextValue 'E8' '  +0.59%  '
Set-TextValue 'E9' '  -0.44%  '
Set-TextValue 'E10' '  -0.54%  '
Set-TextValue 'D11' '5.88'
Set-TextValue 'E11' '  +2.83%  '
Set-TextValue 'E12' '  -1.50%  '
Set-TextValue 'E13' '  -2.50%  '
Set-TextValue 'D14' '37.24'
Set-TextValue 'E14' '  +0.35%  '
Set-TextValue 'D15' '3.661.74'
Set-TextValue 'E15' '  +0.72%  '
Set-TextValue 'E16' '  -1.27%  '
Set-TextValue 'D17' '7.39'
Set-TextValue 'E17' '  +3.34%  '
Set-TextValue 'D18' '63.734.59'
Set-TextValue 'E18' '  +0.12%  '
Set-TextValue 'D19' '3.137.30'
Set-TextValue 'E19' '  +0.70%  '
Set-TextValue 'D20' '467.75'
Set-TextValue 'E20' '  +0.44%  '
Set-TextValue 'E21' '  +0.23%  '
Set-TextValue 'E22' '  +0.15%  '
Set-TextValue 'D23' '7.53'
Set-TextValue 'E23' '  -0.07%  '
Set-TextValue 'B24' 'InternetComputer(DFINITY)'
Set-TextValue 'C24' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D24' '13.00'
Set-TextValue 'E24' '  -1.31%  '
Set-TextValue 'B25' 'Litecoin'
Set-TextValue 'C25' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D25' '81.61'
Set-TextValue 'E25' '  -0.62%  '
Set-TextValue 'E26' '  +6.97%  '
Set-TextValue 'E27' '  +0.02%  '
Set-TextValue 'D28' '9.74'
Set-TextValue 'E28' '  +8.58%  '
Set-TextValue 'D29' '7.46'
Set-TextValue 'E29' '  +9.05%  '
Set-TextValue 'E30' '  +0.15%  '
Set-TextValue 'D31' '2.24'
Set-TextValue 'E31' '  +0.39%  '
Set-TextValue 'E32' '  +0.20%  '
Set-TextValue 'D33' '27.85'
Set-TextValue 'E33' '  +3.11%  '
Set-TextValue 'E34' '  +0.40%  '
Set-TextValue 'E35' '  -4.05%  '
Set-TextValue 'E36' '  +1.27%  '
Set-TextValue 'D37' '6.16'
Set-TextValue 'E37' '  +1.38%  '
Set-TextValue 'E38' '  -3.24%  '
Set-TextValue 'D39' '3.20'
Set-TextValue 'E39' '  -6.69%  '
Set-TextValue 'B40' 'Cosmos'
Set-TextValue 'C40' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D40' '9.40'
Set-TextValue 'E40' '  +8.12%  '
Set-TextValue 'B41' 'OKB'
Set-TextValue 'C41' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D41' '51.56'
Set-TextValue 'E41' '  +1.36%  '
Set-TextValue 'D42' '453.88'
Set-TextValue 'E42' '  -0.01%  '
Set-TextValue 'E43' '  +5.19%  '
Set-TextValue 'D44' '0.0372'
Set-TextValue 'E44' '  -0.06%  '
Set-TextValue 'D45' '2.911.75'
Set-TextValue 'E45' '  +0.98%  '
Set-TextValue 'D46' '39.56'
Set-TextValue 'E46' '  +10.72%  '
Set-TextValue 'D47' '0.108'
Set-TextValue 'E47' '  -3.14%  '
Set-TextValue 'D48' '132.39'
Set-TextValue 'E48' '  +5.77%  '
Set-TextValue 'E49' '  -0.02%  '
Set-TextValue 'E50' '  +2.90%  '
Set-TextValue 'E51' '  -0.93%  '
